# Automatic data refresh: updates DATA_EXTRACCIO timestamps and the
# associated meteorological readings for 2026-02-16 (07:50 run).
# Note: values that look like bare percentages (e.g. "62%") are written
# with a leading apostrophe so Excel stores them as literal text instead
# of auto-converting them to a numeric percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 07:48:37"
$ws.Range("I2").Value = "4.4 mm"
$ws.Range("O2").Value = "1.3 °C"
$ws.Range("E3").Value = "2026-02-16 07:48:39"
$ws.Range("I3").Value = "1.3 mm"
$ws.Range("L3").Value = "50.4 km/h - 233º 7:14 TU"
$ws.Range("M3").Value = "-0.7 °C 7:29 TU"
$ws.Range("E4").Value = "2026-02-16 07:48:42"
$ws.Range("K4").Value = "0.0 MJ/m2"
$ws.Range("E5").Value = "2026-02-16 07:48:44"
$ws.Range("G5").Value = "135 cm"
$ws.Range("I5").Value = "4.8 mm"
$ws.Range("M5").Value = "-0.4 °C 7:27 TU"
$ws.Range("O5").Value = "-0.9 °C"
$ws.Range("E6").Value = "2026-02-16 07:48:46"
$ws.Range("O6").Value = "6.7 °C"
$ws.Range("E7").Value = "2026-02-16 07:48:49"
$ws.Range("J7").Value = "1014.6 hPa"
$ws.Range("K7").Value = "0.1 MJ/m2"
$ws.Range("E8").Value = "2026-02-16 07:48:51"
$ws.Range("K8").Value = "0.1 MJ/m2"
$ws.Range("E9").Value = "2026-02-16 07:48:54"
$ws.Range("E10").Value = "2026-02-16 07:48:56"
$ws.Range("K10").Value = "0.0 MJ/m2"
$ws.Range("E11").Value = "2026-02-16 07:48:58"
$ws.Range("E12").Value = "2026-02-16 07:49:01"
$ws.Range("E13").Value = "2026-02-16 07:49:03"
$ws.Range("O13").Value = "0.4 °C"
$ws.Range("E14").Value = "2026-02-16 07:49:06"
$ws.Range("H14").Value = "'62%"
$ws.Range("N14").Value = "11.6 °C 7:29 TU"
$ws.Range("E15").Value = "2026-02-16 07:49:08"
$ws.Range("E16").Value = "2026-02-16 07:49:10"
$ws.Range("I16").Value = "2.4 mm"
$ws.Range("E17").Value = "2026-02-16 07:49:13"
$ws.Range("L17").Value = "45.4 km/h - 303º 7:02 TU"
$ws.Range("M17").Value = "6.4 °C 7:29 TU"
$ws.Range("E18").Value = "2026-02-16 07:49:15"
$ws.Range("O18").Value = "3.7 °C"
$ws.Range("E19").Value = "2026-02-16 07:49:18"
$ws.Range("N19").Value = "2.1 °C 7:01 TU"
$ws.Range("E20").Value = "2026-02-16 07:49:20"
$ws.Range("K20").Value = "0.1 MJ/m2"
$ws.Range("E21").Value = "2026-02-16 07:49:23"
$ws.Range("H21").Value = "'84%"
$ws.Range("I21").Value = "0.1 mm"
$ws.Range("N21").Value = "2.3 °C 7:00 TU"
$ws.Range("E22").Value = "2026-02-16 07:49:25"
$ws.Range("K22").Value = "0.1 MJ/m2"
$ws.Range("E23").Value = "2026-02-16 07:49:28"
$ws.Range("I23").Value = "2.6 mm"
$ws.Range("O23").Value = "-0.7 °C"
$ws.Range("E24").Value = "2026-02-16 07:49:30"
$ws.Range("N24").Value = "9.9 °C 7:05 TU"
$ws.Range("E25").Value = "2026-02-16 07:49:33"
$ws.Range("H25").Value = "'77%"
$ws.Range("I25").Value = "0.5 mm"
$ws.Range("K25").Value = "0.1 MJ/m2"
$ws.Range("N25").Value = "-0.6 °C 7:12 TU"
$ws.Range("O25").Value = "0.3 °C"
$ws.Range("E26").Value = "2026-02-16 07:49:35"
$ws.Range("E27").Value = "2026-02-16 07:49:38"
$ws.Range("E28").Value = "2026-02-16 07:49:40"
$ws.Range("J28").Value = "1015.7 hPa"
$ws.Range("E29").Value = "2026-02-16 07:49:42"
$ws.Range("E30").Value = "2026-02-16 07:49:45"
$ws.Range("J30").Value = "1014.2 hPa"
$ws.Range("M30").Value = "8.0 °C 7:12 TU"
$ws.Range("O30").Value = "6.9 °C"
$ws.Range("E31").Value = "2026-02-16 07:49:47"
$ws.Range("H31").Value = "'56%"
$ws.Range("J31").Value = "1012.9 hPa"
$ws.Range("K31").Value = "0.0 MJ/m2"
$ws.Range("E32").Value = "2026-02-16 07:49:50"
$ws.Range("E33").Value = "2026-02-16 07:49:52"
$ws.Range("H33").Value = "'77%"
$ws.Range("O33").Value = "3.8 °C"
$ws.Range("E34").Value = "2026-02-16 07:49:55"
$ws.Range("K34").Value = "0.1 MJ/m2"
$ws.Range("E35").Value = "2026-02-16 07:49:57"
$ws.Range("J35").Value = "1018.8 hPa"
$ws.Range("M35").Value = "7.6 °C 7:26 TU"
$ws.Range("E36").Value = "2026-02-16 07:50:00"
$ws.Range("H36").Value = "'94%"
$ws.Range("E37").Value = "2026-02-16 07:50:02"
$ws.Range("J37").Value = "1017.9 hPa"
$ws.Range("O37").Value = "1.8 °C"
$ws.Range("E38").Value = "2026-02-16 07:50:05"
$ws.Range("K38").Value = "0.0 MJ/m2"
$ws.Range("O38").Value = "5.7 °C"
$ws.Range("E39").Value = "2026-02-16 07:50:07"
$ws.Range("H39").Value = "'75%"
$ws.Range("K39").Value = "0.1 MJ/m2"
$ws.Range("E40").Value = "2026-02-16 07:50:10"
$ws.Range("O40").Value = "2.5 °C"
$ws.Range("E41").Value = "2026-02-16 07:50:12"
$ws.Range("K41").Value = "0.0 MJ/m2"
$ws.Range("E42").Value = "2026-02-16 07:50:15"
$ws.Range("E43").Value = "2026-02-16 07:50:17"
$ws.Range("O43").Value = "2.9 °C"
$ws.Range("E44").Value = "2026-02-16 07:50:19"
$ws.Range("I44").Value = "3.0 mm"
$ws.Range("E45").Value = "2026-02-16 07:50:22"
$ws.Range("I45").Value = "2.9 mm"
$ws.Range("J45").Value = "1019.2 hPa"
$ws.Range("L45").Value = "10.1 km/h - 201º 7:28 TU"
$ws.Range("M45").Value = "3.9 °C 7:29 TU"
$ws.Range("O45").Value = "3.2 °C"
$ws.Range("E46").Value = "2026-02-16 07:50:24"
$ws.Range("K46").Value = "0.0 MJ/m2"
$ws.Range("M46").Value = "13.7 °C 7:29 TU"
$ws.Range("O46").Value = "12.8 °C"
